# Append the 2025-11-16 profit-allocation row to the bottom of the table.
# Source data is a running Date/BTC/KAS allocation log; rows 2-75 already
# hold one row per day through 11/15/2025, so this just continues the
# series with the newly computed split for 11/16/2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 76

# Column A holds the date as literal text (matching every row above it,
# e.g. A75 = "11/15/2025"), not an Excel date serial. A leading apostrophe
# forces text entry instead of date auto-detection, then we reset the
# style back to Normal so no stray "quote prefix" formatting sticks to
# the cell (keeping it identical in shape to the other date cells).
$ws.Range("A" + $newRow).Value = "'11/16/2025"
$ws.Range("A" + $newRow).Style = "Normal"

$ws.Range("B" + $newRow).Value = 0.2029201483094979
$ws.Range("C" + $newRow).Value = 0.7970798516905021
